# Commit: "update database and change read_price algorithm"
#
# The underlying data fix replaces a set of placeholder "-" (dash) cells
# in column K (quarter "فصل چهارم منتهی به 1399/12") with the correct
# computed quantities, and updates a handful of already-numeric "0"
# totals/subtotals in that same column to their correct aggregated
# values, now that the feeding rows are no longer blank/placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose K column changes from the "-" placeholder to a real number,
# or from a stale 0 total to the correct aggregated total.
$kColumnUpdates = @{
    11 = 110028
    12 = 206
    15 = 52839
    16 = 433
    17 = 27302
    18 = 52240
    19 = 43397
    20 = 25687
    21 = 103374
    22 = 415506
    24 = 1862
    28 = 56667
    30 = 40480
    31 = 2226
    32 = 1601
    33 = 41
    34 = 61038
    35 = 163915
    40 = 0
    41 = 579421
    48 = 7655133
    49 = 72614
    52 = 9610066
    53 = 189012
    54 = 3977734
    55 = 12830873
    56 = 12073483
    57 = 9578222
    58 = 9296487
    59 = 65283624
    61 = 209781
    65 = 10399671
    67 = 5170494
    68 = 558419
    69 = 341543
    70 = 8993
    71 = 5169742
    72 = 21858643
    77 = 0
    78 = 87142267
}

foreach ($row in $kColumnUpdates.Keys) {
    $ws.Cells.Item($row, 11).Value = $kColumnUpdates[$row]
}
